$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '26.112.66'
Set-TextCell $ws.Range('E2') '  +3.26%  '
Set-TextCell $ws.Range('D3') '1.598.83'
Set-TextCell $ws.Range('E3') '  +2.08%  '
Set-TextCell $ws.Range('E4') '  -0.04%  '
Set-TextCell $ws.Range('D5') '212.22'
Set-TextCell $ws.Range('E5') '  +2.26%  '
Set-TextCell $ws.Range('E6') '  -0.04%  '
Set-TextCell $ws.Range('D7') '0.484'
Set-TextCell $ws.Range('E7') '  +1.72%  '
Set-TextCell $ws.Range('D8') '0.248'
Set-TextCell $ws.Range('E8') '  +2.28%  '
Set-TextCell $ws.Range('D9') '0.0615'
Set-TextCell $ws.Range('E9') '  +1.47%  '
Set-TextCell $ws.Range('D10') '17.90'
Set-TextCell $ws.Range('E10') '  +0.48%  '
Set-TextCell $ws.Range('E11') '  +4.87%  '
Set-TextCell $ws.Range('D12') '1.822.59'
Set-TextCell $ws.Range('E12') '  +2.23%  '
Set-TextCell $ws.Range('D13') '1.594.32'
Set-TextCell $ws.Range('E13') '  +1.70%  '
Set-TextCell $ws.Range('D14') '4.00'
Set-TextCell $ws.Range('E14') '  -0.57%  '
Set-TextCell $ws.Range('D15') '0.509'
Set-TextCell $ws.Range('E15') '  +0.79%  '
Set-TextCell $ws.Range('D16') '26.106.02'
Set-TextCell $ws.Range('E16') '  +3.28%  '
Set-TextCell $ws.Range('D17') '60.43'
Set-TextCell $ws.Range('E17') '  +1.79%  '
Set-TextCell $ws.Range('D18') '0.0₃0721'
Set-TextCell $ws.Range('E18') '  +1.12%  '
Set-TextCell $ws.Range('D19') '1.00'
Set-TextCell $ws.Range('E19') '  -0.14%  '
Set-TextCell $ws.Range('D20') '204.99'
Set-TextCell $ws.Range('E20') '  +10.34%  '
Set-TextCell $ws.Range('D21') '4.23'
Set-TextCell $ws.Range('E21') '  +2.23%  '
Set-TextCell $ws.Range('D22') '9.30'
Set-TextCell $ws.Range('E22') '  +0.24%  '
Set-TextCell $ws.Range('D23') '5.97'
Set-TextCell $ws.Range('E23') '  +1.63%  '
Set-TextCell $ws.Range('D24') '1.83'
Set-TextCell $ws.Range('E24') '  +11.79%  '
Set-TextCell $ws.Range('E25') '  +1.78%  '
Set-TextCell $ws.Range('D26') '1.01'
Set-TextCell $ws.Range('E26') '  -0.05%  '
Set-TextCell $ws.Range('E27') '  -3.35%  '
Set-TextCell $ws.Range('D28') '15.21'
Set-TextCell $ws.Range('E28') '  +2.27%  '
Set-TextCell $ws.Range('D29') '6.43'
Set-TextCell $ws.Range('E29') '  -0.29%  '
Set-TextCell $ws.Range('E30') '  +1.15%  '
Set-TextCell $ws.Range('D31') '0.0469'
Set-TextCell $ws.Range('E31') '  +1.49%  '
Set-TextCell $ws.Range('D32') '3.12'
Set-TextCell $ws.Range('E32') '  +2.76%  '
Set-TextCell $ws.Range('E33') '  -0.05%  '
Set-TextCell $ws.Range('D34') '1.47'
Set-TextCell $ws.Range('E34') '  +0.76%  '
Set-TextCell $ws.Range('E35') '  +2.48%  '
Set-TextCell $ws.Range('D36') '0.0163'
Set-TextCell $ws.Range('E36') '  +8.90%  '
Set-TextCell $ws.Range('D37') '1.105.88'
Set-TextCell $ws.Range('E37') '  +1.77%  '
Set-TextCell $ws.Range('E39') '  +0.38%  '
Set-TextCell $ws.Range('D40') '0.776'
Set-TextCell $ws.Range('E40') '  +0.45%  '
Set-TextCell $ws.Range('D41') '0.493'
Set-TextCell $ws.Range('E41') '  -0.48%  '
Set-TextCell $ws.Range('E42') '  +1.94%  '
Set-TextCell $ws.Range('D43') '1.736.34'
Set-TextCell $ws.Range('E43') '  +2.32%  '
Set-TextCell $ws.Range('D44') '92.54'
Set-TextCell $ws.Range('E44') '  -0.79%  '
Set-TextCell $ws.Range('D45') '5.08'
Set-TextCell $ws.Range('E45') '  +0.54%  '
Set-TextCell $ws.Range('D46') '0.0₆0105'
Set-TextCell $ws.Range('E46') '  -5.37%  '
Set-TextCell $ws.Range('D47') '1.51'
Set-TextCell $ws.Range('E47') '  +6.35%  '
Set-TextCell $ws.Range('D48') '53.36'
Set-TextCell $ws.Range('E48') '  +1.14%  '
Set-TextCell $ws.Range('E49') '  +0.29%  '
Set-TextCell $ws.Range('E50') '  +0.79%  '
Set-TextCell $ws.Range('E51') '  +0.10%  '
